$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 = Disease Ontology: source_version v2024-08-29 -> v2024-09-27
$ws.Range("E3").Value = "v2024-09-27"

# Row 4 = Experimental Factor Ontology: source_version v3.69.0 -> v3.71.0
$ws.Range("E4").Value = "v3.71.0"

# Update the active selection to E4 (matches the saved selection state in the diff)
$ws.Range("E4").Select()
